# Apply the "Creating pages and test for twitter login" edit.
#
# Summary of changes:
#  - Sheet "test_suite": runmode flags for the existing 3 tests flip from Y to N,
#    and a new row is appended for "TwtLoginTest" with runmode "Y".
#  - Sheet "AddCustomerTest": the two data rows (runmode N/Y) swap their
#    firstname/lastname payload - row2 becomes Rahul/Jadhwani with runmode Y,
#    row3 becomes Deepender/Singh.
#  - Sheet "OpenAccountTest": unchanged.

$wb = $excel.ActiveWorkbook

$wsSuite = $wb.Worksheets.Item("test_suite")
$wsAdd   = $wb.Worksheets.Item("AddCustomerTest")
# "OpenAccountTest" sheet is untouched by this change.

# ---------------------------------------------------------------------------
# test_suite: flip existing runmodes to "N" and append the new TwtLoginTest row
# ---------------------------------------------------------------------------
$wsSuite.Range("B2").Value = "N"
$wsSuite.Range("B3").Value = "N"
$wsSuite.Range("B4").Value = "N"
$wsSuite.Range("A5").Value = "TwtLoginTest"
$wsSuite.Range("B5").Value = "Y"

# ---------------------------------------------------------------------------
# AddCustomerTest: swap the two customer name rows
# ---------------------------------------------------------------------------
$wsAdd.Range("A2").Value = "Y"
$wsAdd.Range("B2").Value = "Rahul"
$wsAdd.Range("C2").Value = "Jadhwani"

$wsAdd.Range("A3").Value = "Y"
$wsAdd.Range("B3").Value = "Deepender"
$wsAdd.Range("C3").Value = "Singh"

# ---------------------------------------------------------------------------
# Selections: mirror the final cursor positions recorded in the workbook.
# Select test_suite first, then finish on AddCustomerTest so it stays the
# active/tabbed sheet (matches the saved workbook state).
# ---------------------------------------------------------------------------
$wsSuite.Range("B5").Select() | Out-Null
$wsAdd.Range("B7").Select() | Out-Null
